# Updates cryptos list values (price + volume(1h) columns) to reflect the
# latest scrape, matching the commit "Updated cryptos list ... with GitHub
# Actions". Also swaps the Bittensor / Fetch.AI rows (32 and 33).
#
# D-column numeric-looking values are written with a leading apostrophe so
# Excel keeps them as literal text (matching the workbook's existing
# plain-text price formatting, e.g. "484.15" instead of being parsed into
# the number 484.15/484.14999999999998). Values that already aren't valid
# numbers (thousand-separated prices like "71.858.64", or the special
# subscript price "0.0₃0940") are left unprefixed since Excel treats them
# as text automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.858.64"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "2.619.05"
$ws.Range("E3").Value = "  -2.18%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'585.40"
$ws.Range("E5").Value = "  -1.67%  "
$ws.Range("D6").Value = "'173.69"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D8").Value = "'0.517"
$ws.Range("E8").Value = "  -0.55%  "
$ws.Range("D9").Value = "2.618.84"
$ws.Range("E9").Value = "  -2.22%  "
$ws.Range("D10").Value = "'0.169"
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "'0.357"
$ws.Range("E12").Value = "  +1.74%  "
$ws.Range("D13").Value = "'4.89"
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("D14").Value = "'0.0000189"
$ws.Range("E14").Value = "  +3.17%  "
$ws.Range("D15").Value = "3.099.06"
$ws.Range("D16").Value = "71.718.60"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").Value = "'25.59"
$ws.Range("E17").Value = "  -1.74%  "
$ws.Range("D18").Value = "2.605.72"
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("D19").Value = "'12.04"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("D21").Value = "'373.68"
$ws.Range("E21").Value = "  +1.59%  "
$ws.Range("E22").Value = "  -1.70%  "
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("D24").Value = "'71.27"
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("E26").Value = "  -2.35%  "
$ws.Range("D27").Value = "'9.30"
$ws.Range("E27").Value = "  -5.14%  "
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "0.0₃0940"
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("D31").Value = "'7.93"
$ws.Range("E31").Value = "  -1.25%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'484.15"
$ws.Range("E32").Value = "  -3.66%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.31"
$ws.Range("E33").Value = "  +2.42%  "
$ws.Range("E34").Value = "  -0.79%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'159.40"
$ws.Range("E36").Value = "  -2.01%  "
$ws.Range("E37").Value = "  +7.39%  "
$ws.Range("D38").Value = "'19.22"
$ws.Range("E38").Value = "  -0.71%  "
$ws.Range("D39").Value = "'18.91"
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("E40").Value = "  -1.16%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E42").Value = "  -4.30%  "
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("D44").Value = "'4.83"
$ws.Range("E44").Value = "  -2.88%  "
$ws.Range("E45").Value = "  -1.72%  "
$ws.Range("D46").Value = "'38.95"
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("D47").Value = "'149.43"
$ws.Range("E47").Value = "  -4.02%  "
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("E49").Value = "  -1.30%  "
$ws.Range("E50").Value = "  -3.89%  "
$ws.Range("D51").Value = "'0.602"
$ws.Range("E51").Value = "  -0.01%  "

